$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row to the new short column codes.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Title-case the Spanish connector words ("de", "del", "la", "el",
#    "los", "las", "y") wherever they appear *lowercase* inside the
#    state/municipality name strings (columns A and B), for every data
#    row. Words that are already capitalized (e.g. a leading "La Paz")
#    must stay untouched, so match on exact (case-sensitive) lowercase
#    tokens only.
$connectors = @('de', 'del', 'la', 'el', 'los', 'las', 'y')
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @(1, 2)) {
        $cell = $ws.Cells.Item($r, $col)
        $v = $cell.Value()
        if ($v -ne $null -and $v.GetType().Name -eq "String") {
            $words = $v -split ' '
            $changed = $false
            $out = @()
            foreach ($w in $words) {
                if ($connectors -contains $w) {
                    $out += ($w.Substring(0, 1).ToUpper() + $w.Substring(1))
                    $changed = $true
                } else {
                    $out += $w
                }
            }
            if ($changed) {
                $cell.Value = ($out -join ' ')
            }
        }
    }
}

# 3) Tiny floating-point re-round picked up by the source pipeline.
$ws.Range("D1172").Value = 0.009166107757657051

# 4) Drop the trailing metadata/footer rows (sample size, source,
#    author, date) that used to sit below the data table.
$ws.Range("A1226:A1230").EntireRow.Delete()
